$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original (pre-edit) data for columns D,K,L,M,N,O,P,Q,S,T for rows 156..230, captured as literal
# values (these are the values that exist in the workbook before this script runs).
# Index 0 corresponds to row 156, index 74 corresponds to row 230.
$src = @(
    @{D=44249; K="Valencia"; L="Primera"; M=100; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44265; K="Valencia"; L="Primera"; M=160; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44343; K="Fukumoto"; L="Primera"; M=400; N=17000; O=18000; P=17500; Q="`$/caja 15 kilos empedrada"; S=1167; T=15}
    @{D=44343; K="Fukumoto"; L="Segunda"; M=200; N=14000; O=14000; P=14000; Q="`$/caja 15 kilos empedrada"; S=933; T=15}
    @{D=44201; K="Valencia"; L="Primera"; M=400; N=20000; O=21000; P=20500; Q="`$/caja 15 kilos empedrada"; S=1367; T=15}
    @{D=44280; K="Valencia"; L="Primera"; M=160; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44447; K="Navel Late"; L="Primera"; M=200; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44447; K="Navel Late"; L="Segunda"; M=100; N=10000; O=10000; P=10000; Q="`$/caja 15 kilos empedrada"; S=667; T=15}
    @{D=44270; K="Valencia"; L="Primera"; M=120; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44267; K="Valencia"; L="Primera"; M=200; N=24000; O=24000; P=24000; Q="`$/caja 15 kilos empedrada"; S=1600; T=15}
    @{D=44187; K="Valencia"; L="Primera"; M=600; N=20000; O=21000; P=20500; Q="`$/caja 15 kilos empedrada"; S=1367; T=15}
    @{D=44386; K="New Hall"; L="Primera"; M=400; N=12000; O=13000; P=12500; Q="`$/caja 15 kilos empedrada"; S=833; T=15}
    @{D=44386; K="New Hall"; L="Segunda"; M=200; N=10000; O=10000; P=10000; Q="`$/caja 15 kilos empedrada"; S=667; T=15}
    @{D=44386; K="New Hall"; L="Segunda"; M=400; N=9000; O=9000; P=9000; Q="`$/malla 18 kilos"; S=500; T=18}
    @{D=44308; K="Valencia"; L="Primera"; M=400; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44264; K="Valencia"; L="Primera"; M=200; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44463; K="Lane Late"; L="Primera"; M=500; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44463; K="Lane Late"; L="Segunda"; M=250; N=10000; O=10000; P=10000; Q="`$/caja 15 kilos empedrada"; S=667; T=15}
    @{D=44196; K="Valencia"; L="Primera"; M=200; N=20000; O=21000; P=20500; Q="`$/caja 15 kilos empedrada"; S=1367; T=15}
    @{D=44301; K="Valencia"; L="Primera"; M=200; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44251; K="Valencia"; L="Primera"; M=160; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44243; K="Valencia"; L="Primera"; M=400; N=24000; O=25000; P=24500; Q="`$/malla 16 kilos"; S=1531; T=16}
    @{D=44252; K="Valencia"; L="Primera"; M=120; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44166; K="Valencia"; L="Primera"; M=600; N=16000; O=17000; P=16500; Q="`$/caja 15 kilos empedrada"; S=1100; T=15}
    @{D=44168; K="Valencia"; L="Primera"; M=600; N=16500; O=17000; P=16750; Q="`$/caja 15 kilos empedrada"; S=1117; T=15}
    @{D=44369; K="New Hall"; L="Primera"; M=600; N=14000; O=15000; P=14500; Q="`$/caja 15 kilos empedrada"; S=967; T=15}
    @{D=44369; K="New Hall"; L="Segunda"; M=300; N=11000; O=11000; P=11000; Q="`$/caja 15 kilos empedrada"; S=733; T=15}
    @{D=44369; K="New Hall"; L="Segunda"; M=500; N=10000; O=10000; P=10000; Q="`$/malla 18 kilos"; S=556; T=18}
    @{D=44433; K="New Hall"; L="Primera"; M=120; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44433; K="New Hall"; L="Segunda"; M=100; N=10000; O=10000; P=10000; Q="`$/caja 15 kilos empedrada"; S=667; T=15}
    @{D=44316; K="Valencia"; L="Primera"; M=400; N=23000; O=24000; P=23500; Q="`$/caja 15 kilos empedrada"; S=1567; T=15}
    @{D=44316; K="Valencia"; L="Segunda"; M=100; N=20000; O=20000; P=20000; Q="`$/caja 15 kilos empedrada"; S=1333; T=15}
    @{D=44279; K="Valencia"; L="Primera"; M=160; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44397; K="New Hall"; L="Segunda"; M=200; N=10000; O=10000; P=10000; Q="`$/caja 15 kilos"; S=667; T=15}
    @{D=44397; K="New Hall"; L="Segunda"; M=400; N=9000; O=9000; P=9000; Q="`$/malla 18 kilos"; S=500; T=18}
    @{D=44363; K="New Hall"; L="Primera"; M=200; N=15000; O=16000; P=15500; Q="`$/caja 15 kilos empedrada"; S=1033; T=15}
    @{D=44277; K="Valencia"; L="Primera"; M=200; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44273; K="Valencia"; L="Primera"; M=200; N=25000; O=25000; P=25000; Q="`$/caja 15 kilos empedrada"; S=1667; T=15}
    @{D=44438; K="Navel Late"; L="Primera"; M=200; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44372; K="New Hall"; L="Primera"; M=400; N=14000; O=15000; P=14500; Q="`$/caja 15 kilos empedrada"; S=967; T=15}
    @{D=44372; K="New Hall"; L="Segunda"; M=200; N=11000; O=11000; P=11000; Q="`$/caja 15 kilos empedrada"; S=733; T=15}
    @{D=44372; K="New Hall"; L="Segunda"; M=400; N=10000; O=10000; P=10000; Q="`$/malla 18 kilos"; S=556; T=18}
    @{D=44286; K="Valencia"; L="Primera"; M=160; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44209; K="Valencia"; L="Primera"; M=160; N=22000; O=23000; P=22500; Q="`$/caja 15 kilos empedrada"; S=1500; T=15}
    @{D=44160; K="Valencia"; L="Primera"; M=100; N=17000; O=17000; P=17000; Q="`$/caja 15 kilos empedrada"; S=1133; T=15}
    @{D=44160; K="Valencia"; L="Segunda"; M=100; N=16500; O=16500; P=16500; Q="`$/caja 15 kilos empedrada"; S=1100; T=15}
    @{D=44351; K="Fukumoto"; L="Primera"; M=600; N=15500; O=16000; P=15750; Q="`$/caja 15 kilos empedrada"; S=1050; T=15}
    @{D=44351; K="Fukumoto"; L="Segunda"; M=200; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44365; K="New Hall"; L="Primera"; M=600; N=15000; O=16000; P=15500; Q="`$/caja 15 kilos empedrada"; S=1033; T=15}
    @{D=44365; K="New Hall"; L="Segunda"; M=300; N=12500; O=12500; P=12500; Q="`$/caja 15 kilos empedrada"; S=833; T=15}
    @{D=44306; K="Valencia"; L="Primera"; M=400; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44215; K="Valencia"; L="Primera"; M=200; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44175; K="Valencia"; L="Primera"; M=200; N=16500; O=17000; P=16750; Q="`$/caja 15 kilos empedrada"; S=1117; T=15}
    @{D=44357; K="Fukumoto"; L="Primera"; M=200; N=15000; O=15500; P=15250; Q="`$/caja 15 kilos empedrada"; S=1017; T=15}
    @{D=44357; K="Fukumoto"; L="Segunda"; M=100; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44203; K="Valencia"; L="Primera"; M=300; N=20000; O=21000; P=20500; Q="`$/caja 15 kilos empedrada"; S=1367; T=15}
    @{D=44162; K="Valencia"; L="Primera"; M=400; N=17000; O=17500; P=17250; Q="`$/caja 15 kilos empedrada"; S=1150; T=15}
    @{D=44411; K="New Hall"; L="Primera"; M=400; N=14000; O=14000; P=14000; Q="`$/caja 15 kilos empedrada"; S=933; T=15}
    @{D=44411; K="New Hall"; L="Segunda"; M=200; N=10000; O=10000; P=10000; Q="`$/caja 15 kilos empedrada"; S=667; T=15}
    @{D=44411; K="New Hall"; L="Segunda"; M=400; N=8000; O=8000; P=8000; Q="`$/malla 18 kilos"; S=444; T=18}
    @{D=44257; K="Valencia"; L="Primera"; M=300; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44244; K="Valencia"; L="Primera"; M=120; N=24000; O=25000; P=24500; Q="`$/caja 15 kilos empedrada"; S=1633; T=15}
    @{D=44176; K="Valencia"; L="Primera"; M=400; N=16500; O=17000; P=16750; Q="`$/caja 15 kilos empedrada"; S=1117; T=15}
    @{D=44239; K="Valencia"; L="Primera"; M=400; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44376; K="New Hall"; L="Primera"; M=400; N=14000; O=15000; P=14500; Q="`$/caja 15 kilos empedrada"; S=967; T=15}
    @{D=44376; K="New Hall"; L="Primera"; M=500; N=10000; O=10000; P=10000; Q="`$/malla 18 kilos"; S=556; T=18}
    @{D=44376; K="New Hall"; L="Segunda"; M=200; N=11000; O=11000; P=11000; Q="`$/caja 15 kilos empedrada"; S=733; T=15}
    @{D=44292; K="Valencia"; L="Primera"; M=400; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44358; K="New Hall"; L="Primera"; M=600; N=16000; O=17000; P=16500; Q="`$/caja 15 kilos empedrada"; S=1100; T=15}
    @{D=44358; K="New Hall"; L="Segunda"; M=300; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44211; K="Valencia"; L="Primera"; M=300; N=25000; O=26000; P=25500; Q="`$/caja 15 kilos empedrada"; S=1700; T=15}
    @{D=44425; K="New Hall"; L="Primera"; M=400; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
    @{D=44425; K="New Hall"; L="Segunda"; M=200; N=9000; O=9000; P=9000; Q="`$/caja 15 kilos empedrada"; S=600; T=15}
    @{D=44425; K="New Hall"; L="Segunda"; M=400; N=7000; O=7000; P=7000; Q="`$/malla 18 kilos"; S=389; T=18}
    @{D=44323; K="Fukumoto"; L="Primera"; M=600; N=19000; O=20000; P=19500; Q="`$/caja 15 kilos empedrada"; S=1300; T=15}
)

# New data to be inserted for the two brand-new report rows (156 and 157). All existing rows
# from 158..232 are simply the old rows 156..230 shifted down by two positions.
$newRow156 = @{D=44466; K="Lane Late"; L="Primera"; M=200; N=13000; O=13000; P=13000; Q="`$/caja 15 kilos empedrada"; S=867; T=15}
$newRow157 = @{D=44466; K="Lane Late"; L="Segunda"; M=100; N=10000; O=10000; P=10000; Q="`$/caja 15 kilos empedrada"; S=667; T=15}

function Set-DataRow($rowNum, $data) {
    $ws.Range("D$rowNum").Value = $data.D
    $ws.Range("D$rowNum").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Range("K$rowNum").Value = $data.K
    $ws.Range("L$rowNum").Value = $data.L
    $ws.Range("M$rowNum").Value = $data.M
    $ws.Range("N$rowNum").Value = $data.N
    $ws.Range("O$rowNum").Value = $data.O
    $ws.Range("P$rowNum").Value = $data.P
    $ws.Range("Q$rowNum").Value = $data.Q
    $ws.Range("S$rowNum").Value = $data.S
    $ws.Range("T$rowNum").Value = $data.T
}

# Columns that stay constant across every data row in this sheet.
$constA = 4
$constB = "Feria Lagunitas de Puerto Montt"
$constC = "Los Lagos"
$constE = 10
$constF = "Fruta"
$constG = 100102
$constH = "Cítricos"
$constI = 100102005
$constJ = "Naranja"
$constR = "Región de O'Higgins"

function Set-ConstCols($rowNum) {
    $ws.Range("A$rowNum").Value = $constA
    $ws.Range("B$rowNum").Value = $constB
    $ws.Range("C$rowNum").Value = $constC
    $ws.Range("E$rowNum").Value = $constE
    $ws.Range("F$rowNum").Value = $constF
    $ws.Range("G$rowNum").Value = $constG
    $ws.Range("H$rowNum").Value = $constH
    $ws.Range("I$rowNum").Value = $constI
    $ws.Range("J$rowNum").Value = $constJ
    $ws.Range("R$rowNum").Value = $constR
}

# Two new rows are appended at the bottom of the sheet (231 and 232); give them the constant
# columns explicitly since they do not exist yet.
Set-ConstCols 231
Set-ConstCols 232

# Shift every existing data row (156..230) down by two rows. Walk from the bottom (232) up to
# 158 so that each row's source data (row N-2) is read before it is overwritten.
for ($n = 232; $n -ge 158; $n--) {
    $data = $src[$n - 2 - 156]
    Set-DataRow $n $data
}

# Rows 156 and 157 become brand-new report rows with fresh data.
Set-DataRow 156 $newRow156
Set-DataRow 157 $newRow157
